$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # Locate the target text with Find, then assign .Text directly on the
    # found range rather than passing a replacement string into
    # Find.Execute -- the latter runs the replacement through Word's
    # "replace as you type" autoformatting (smart quotes, etc.) which we
    # do not want here.
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

# Mission / Vision paragraph: "not for & in" -> "not for, and in"
Replace-Text "We are former and current government workers, designers, technologists, and concerned citizens who believe in building with, not for & in partnering directly with users both in and outside of government to find solutions that help bridge the gap between the two." "We are former and current government workers, designers, technologists, and concerned citizens who believe in building with, not for, and in partnering directly with users both in and outside of government to find solutions that help bridge the gap between the two."

# Free Law Founders paragraph: "us & the" -> "us and the"; " -- to" -> em dash "to"
Replace-Text "We're building a national network of people like us & the Free Law Founders -- to share resources and expertise around solving shared challenges in modernizing the lawmaking process." ("We're building a national network of people like us and the Free Law Founders" + [char]0x2014 + "to share resources and expertise around solving shared challenges in modernizing the lawmaking process.")

# Reimbursement intro paragraph: " - including" / "subscriptions, must" -> em dashes
Replace-Text "The OpenGov Foundation will reimburse pre-approved, work-related expenses. Any expenditure totalling more than `$1,000 - including cumulative expenses like software subscriptions, must be pre-approved by the Board of Directors and Executive Director. Sub-`$1,000 reimbursable expenses which do not require BOD approval include but are not limited to:" ("The OpenGov Foundation will reimburse pre-approved, work-related expenses. Any expenditure totalling more than `$1,000" + [char]0x2014 + "including cumulative expenses like software subscriptions" + [char]0x2014 + "must be pre-approved by the Board of Directors and Executive Director. Sub-`$1,000 reimbursable expenses which do not require BOD approval include but are not limited to:")

# Expense list items -> add trailing period
Replace-Text "Travel between work meetings and the office" "Travel between work meetings and the office."
Replace-Text "Travel to, from, and during conferences and events" "Travel to, from, and during conferences and events."
Replace-Text "Housing and food during conferences and events" "Housing and food during conferences and events."
Replace-Text "Business-meeting meals and drinks" "Business-meeting meals and drinks."
Replace-Text "Office supplies" "Office supplies."
Replace-Text "Organization-specific subscriptions and/or software" "Organization-specific subscriptions and/or software."
Replace-Text "Expenses stemming from events organized by The OpenGov Foundation" "Expenses stemming from events organized by The OpenGov Foundation."

# Travel policy sentence: "reasonably priced" -> "reasonably-priced"
Replace-Text "Employees are expected to travel coach class when traveling by air, book reasonably priced hotels and not spend excessively on meals and other travel expenses." "Employees are expected to travel coach class when traveling by air, book reasonably-priced hotels and not spend excessively on meals and other travel expenses."

# Work hours paragraph: "11am-4pm" -> "11 AM" en-dash "4 PM"
Replace-Text "Employees should check in on Teamwork by 11 AM EST each week day with an update of their work plan for the day. Communication is vital to the success of our team; therefore, team members are expected to be reasonably accessible by phone, Slack, and email between the hours of 11am-4pm." ("Employees should check in on Teamwork by 11 AM EST each week day with an update of their work plan for the day. Communication is vital to the success of our team; therefore, team members are expected to be reasonably accessible by phone, Slack, and email between the hours of 11 AM" + [char]0x2013 + "4 PM.")

# Office hours paragraph: "(9am - 5pm) - for ... etc - should" -> "(9 PM" en-dash "5 PM)" em-dash "for ... etc." em-dash "should"
Replace-Text "Any time that you won't be available during usual office hours (9am - 5pm) - for appointments, sick time, vacation, etc - should be put on the calendar, and an email must be sent to the employee's supervisor. For non-emergency medical absences, please notify your supervisor at least 48 hours in advance." ("Any time that you won't be available during usual office hours (9 PM" + [char]0x2013 + "5 PM)" + [char]0x2014 + "for appointments, sick time, vacation, etc." + [char]0x2014 + "should be put on the calendar, and an email must be sent to the employee's supervisor. For non-emergency medical absences, please notify your supervisor at least 48 hours in advance.")

# Gmail / Google+ paragraph
Replace-Text "We will set you up with an opengovfoundation.org gmail account. This will provide you access to our email, calendar, and Google Plus accounts. Please use this for all work-related correspondence, scheduling, etc. Personal accounts should not be used for work-related activities." "We will set you up with an opengovfoundation.org Gmail account. This will provide you access to our email, calendar, and Google+ accounts. Please use this for all work-related correspondence, scheduling, etc. Personal accounts should not be used for work-related activities."
